# CAMERA_rules_neg.xlsx edit:
#  1. Fix copy/paste bug: four "loss of hexose fragments" rules were
#     mistakenly annotated as positive-mode ions ([M+H-...]+) even though
#     this rule set is for negative mode. Re-label them as negative mode.
#  2. Add a new rule for the loss of two hexoses: [M-H-(Hexose-H2O)2]-

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. fix the mis-annotated rows (140-143) ---
$ws.Range("A140").Value2 = "[M-H-Hexose-H2O-CH4]-"
$ws.Range("A141").Value2 = "[M-H-Hexose-H2O-CH3OH]-"
$ws.Range("A142").Value2 = "[M-H-Hexose-H2O-C2H4O]-"
$ws.Range("A143").Value2 = "[M-H-Hexose-H2O-C2H4O2]-"

# --- 2. append the new rule in row 144 ---
# copy the formatting of the row above down into the new row first
$ws.Range("A143:G143").Copy()
$ws.Range("A144:G144").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A144").Value2 = "[M-H-(Hexose-H2O)2]-"
$ws.Range("B144").Value2 = 1
$ws.Range("C144").Value2 = 1
$ws.Range("D144").Formula = "=-(2*162.05283-D2)"
$ws.Range("E144").Value2 = 178
$ws.Range("F144").Value2 = 0
$ws.Range("G144").Value2 = 0.5

# --- cosmetic: move the window selection down to the newly added rows,
#     matching the author's view state when they made the edit ---
$ws.Range("A147").Select() | Out-Null

$excel.Calculate() | Out-Null
